$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "BISCUT"
$ws.Range("B3").Value = "TOFEE"
$ws.Range("B2").Value = "CAKE"

$ws.Range("D10").Select()
